# Insert a new data row at row 96 (pushing existing rows 96-197 down to 97-198)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 96 - this shifts rows 96:197 down to 97:198.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new record.
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C96").Value = "Ñuble"
$ws.Range("D96").Value = 44587
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 100112043
$ws.Range("G96").Value = "Pepino ensalada"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 120
$ws.Range("K96").Value = 8000
$ws.Range("L96").Value = 8500
$ws.Range("M96").Value = 8250
$ws.Range("N96").Value = "$/caja 80 unidades"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 103
$ws.Range("Q96").Value = 80
$ws.Range("R96").Value = "Hortaliza"
